# "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# The "Periodo Mora" / "Valor Mora" table (rows 16-22, columns E/F) was
# refreshed and is now listed in ascending period order instead of the
# previous descending order. The values themselves are unchanged -
# row-for-row they are simply reversed (top<->bottom), which also swaps
# which row carries the one-off "Valor Mora" amount of 25439 (the rest
# are 36341).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodo Mora column (E16:E22): reverse descending -> ascending order.
$ws.Range("E16").Value = "2108"
$ws.Range("E17").Value = "2109"
$ws.Range("E18").Value = "2110"
$ws.Range("E19").Value = "2111"
$ws.Range("E20").Value = "2112"
$ws.Range("E21").Value = "2201"
$ws.Range("E22").Value = "2202"

# Valor Mora column (F16:F22): the odd-one-out amount moves from the
# first data row to the last, following its period (2202).
$ws.Range("F16").Value = 36341
$ws.Range("F22").Value = 25439
